$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new header columns were added to the experiment type table: "Pattern"
# and "Pattern Type", following the existing "Unique" column (L4).
# Copy L4's formatting (bold header style) onto the new header cells, then
# set their text.
$ws.Range("L4").Copy()
$ws.Range("M4:N4").PasteSpecial(-4122)

$ws.Range("M4").Value = "Pattern"
$ws.Range("N4").Value = "Pattern Type"

# Mirror the cursor/selection ending up on the newly added header cells,
# as it does in the edited workbook.
$ws.Range("M4:N4").Select()
